$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1439
$ws.Range("I98").Value = 1507.9546
$ws.Range("K98").Value = 1507.9546
$ws.Range("M98").Value = -9.954600000000028

$ws.Range("H122").Value = 1439
$ws.Range("I122").Value = 1507.9546
$ws.Range("K122").Value = 4523.8638
$ws.Range("M122").Value = -2073.8638

$ws.Range("H137").Value = 3212.5
$ws.Range("I137").Value = 1425
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 4275
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -1725
$ws.Range("N137").Value = -20100

$ws.Range("H138").Value = 7980.968
$ws.Range("I138").Value = 4323.778
$ws.Range("J138").Value = 9477.091
$ws.Range("K138").Value = 12971.334
$ws.Range("L138").Value = 28431.273
$ws.Range("M138").Value = -7831.334000000001
$ws.Range("N138").Value = -38711.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6361.7207
$ws.Range("I32").Value = 6361.7207
$ws.Range("K32").Value = 6361.7207
$ws.Range("M32").Value = -6074.7207

$ws.Range("H61").Value = 55558716
$ws.Range("I61").Value = 55558716
$ws.Range("K61").Value = 55558716
$ws.Range("M61").Value = -55558504

$ws.Range("H74").Value = 5927.9663
$ws.Range("I74").Value = 5927.9663
$ws.Range("K74").Value = 5927.9663
$ws.Range("M74").Value = -5053.9663

$ws.Range("H77").Value = 5927.9663
$ws.Range("I77").Value = 5927.9663
$ws.Range("K77").Value = 29639.8315
$ws.Range("M77").Value = -25271.8315

$ws.Range("H102").Value = 1685228.1
$ws.Range("I102").Value = 2180284.5
$ws.Range("J102").Value = 2036.6
$ws.Range("K102").Value = 2180284.5
$ws.Range("L102").Value = 2036.6
$ws.Range("M102").Value = -2178662.5
$ws.Range("N102").Value = -5280.6

$ws.Range("H136").Value = 55558716
$ws.Range("I136").Value = 55558716
$ws.Range("K136").Value = 166676148
$ws.Range("M136").Value = -166673598

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3671.7144
$ws.Range("I20").Value = 2193.2856
$ws.Range("J20").Value = 5150.143
$ws.Range("K20").Value = 2193.2856
$ws.Range("L20").Value = 5150.143
$ws.Range("M20").Value = -1946.2856
$ws.Range("N20").Value = -5644.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12202.56
$ws.Range("I31").Value = 8719.223
$ws.Range("J31").Value = 14161.9375
$ws.Range("K31").Value = 8719.223
$ws.Range("L31").Value = 14161.9375
$ws.Range("M31").Value = -8424.223
$ws.Range("N31").Value = -14751.9375

$ws.Range("H34").Value = 12202.56
$ws.Range("I34").Value = 8719.223
$ws.Range("J34").Value = 14161.9375
$ws.Range("K34").Value = 8719.223
$ws.Range("L34").Value = 14161.9375
$ws.Range("M34").Value = -8517.223
$ws.Range("N34").Value = -14565.9375

$ws.Range("H107").Value = 866123.2
$ws.Range("I107").Value = 1276310
$ws.Range("J107").Value = 45749.5
$ws.Range("K107").Value = 1276310
$ws.Range("L107").Value = 45749.5
$ws.Range("M107").Value = -1274390
$ws.Range("N107").Value = -49589.5

$ws.Range("H132").Value = 2449.7932
$ws.Range("I132").Value = 2249.88
$ws.Range("J132").Value = 3699.25
$ws.Range("K132").Value = 6749.64
$ws.Range("L132").Value = 11097.75
$ws.Range("M132").Value = -4219.64
$ws.Range("N132").Value = -16157.75

$ws.Range("H141").Value = 500498.34
$ws.Range("J141").Value = 525559.75
$ws.Range("L141").Value = 525559.75
$ws.Range("N141").Value = -535919.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 9095589
$ws.Range("J137").Value = 3653
$ws.Range("L137").Value = 10959
$ws.Range("N137").Value = -21159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2679.8
$ws.Range("I80").Value = 2833.3333
$ws.Range("J80").Value = 2449.5
$ws.Range("K80").Value = 2833.3333
$ws.Range("L80").Value = 2449.5
$ws.Range("M80").Value = -1835.3333
$ws.Range("N80").Value = -4445.5

$ws.Range("H83").Value = 2679.8
$ws.Range("I83").Value = 2833.3333
$ws.Range("J83").Value = 2449.5
$ws.Range("K83").Value = 14166.6665
$ws.Range("L83").Value = 12247.5
$ws.Range("M83").Value = -9174.666499999999
$ws.Range("N83").Value = -22231.5

$ws.Range("H97").Value = 1298
$ws.Range("I97").Value = 996.6667
$ws.Range("K97").Value = 996.6667
$ws.Range("M97").Value = -500.6667

$ws.Range("H122").Value = 102470.914
$ws.Range("I122").Value = 102470.914
$ws.Range("K122").Value = 307412.742
$ws.Range("M122").Value = -304962.742

$ws.Range("H132").Value = 8009
$ws.Range("I132").Value = 7680.6665
$ws.Range("J132").Value = 8747.75
$ws.Range("K132").Value = 23041.9995
$ws.Range("L132").Value = 26243.25
$ws.Range("M132").Value = -20511.9995
$ws.Range("N132").Value = -31303.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2738.3333
$ws.Range("I16").Value = 694.2857
$ws.Range("K16").Value = 694.2857
$ws.Range("M16").Value = -524.2857

$ws.Range("H42").Value = 29899.5
$ws.Range("I42").Value = 29899.5
$ws.Range("K42").Value = 29899.5
$ws.Range("M42").Value = -29336.5

$ws.Range("H49").Value = 29899.5
$ws.Range("I49").Value = 29899.5
$ws.Range("K49").Value = 29899.5
$ws.Range("M49").Value = -29752.5

$ws.Range("H55").Value = 550.6667
$ws.Range("I55").Value = 402.47058
$ws.Range("J55").Value = 683.2632
$ws.Range("K55").Value = 402.47058
$ws.Range("L55").Value = 683.2632
$ws.Range("M55").Value = -229.47058
$ws.Range("N55").Value = -1029.2632

$ws.Range("H82").Value = 2758.5715
$ws.Range("I82").Value = 2495
$ws.Range("K82").Value = 2495
$ws.Range("M82").Value = -2134

$ws.Range("H85").Value = 2758.5715
$ws.Range("I85").Value = 2495
$ws.Range("K85").Value = 2495
$ws.Range("M85").Value = -1247

$ws.Range("H122").Value = 4953.2925
$ws.Range("I122").Value = 4646.2856
$ws.Range("K122").Value = 13938.8568
$ws.Range("M122").Value = -11488.8568

$ws.Range("H132").Value = 65745716
$ws.Range("I132").Value = 76703000
$ws.Range("J132").Value = 2005
$ws.Range("K132").Value = 230109000
$ws.Range("L132").Value = 6015
$ws.Range("M132").Value = -230106470
$ws.Range("N132").Value = -11075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3002.3076
$ws.Range("I122").Value = 1522.3158
$ws.Range("J122").Value = 7019.4287
$ws.Range("K122").Value = 4566.9474
$ws.Range("L122").Value = 21058.2861
$ws.Range("M122").Value = -2116.9474
$ws.Range("N122").Value = -25958.2861
